# CS 536 spring adaptation: bump the cached date-placeholder text on every
# master/layout/notes-master, and update the title slide's lecture number.

$p = $ppt.ActivePresentation
$newDate = "12/20/2024"

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.Name -like "Date Placeholder*") {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DateShape $master.Shapes

# Every slide layout hanging off the master.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShape $layout.Shapes
}

# Notes master's date placeholder.
$notesMaster = $p.NotesMaster
Update-DateShape $notesMaster.Shapes

# Title slide: bump the lecture number for the spring term.
$slide1 = $p.Slides.Item(1)
$titleShape = $slide1.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "22,23 – Gender + Disability"
